# Update the TA with the new world names
# The "ipaddress" row's value (B2) changes from the old simulation-world
# IP to the new one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "192.168.1.131"

# Leave the selection on B5, matching where the author's cursor ended up.
$ws.Range("B5").Select()
